$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for "Date and Time" (pushes "Total time taken for
# the ride" and everything below down by one row, preserving the time style).
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 16:37:01.599000 to 2024-03-11 17:36:23.668000"

# Insert a new row before "Time spent in 0-10 km/h" (currently row 36) for the
# new "Cycle Count of battery" parameter.
$ws.Rows.Item(35).Insert()

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 42
